# Actualización automática hashcode jue sep  3 02:10:06 CEST 2020
# Refresh the hashcode column (B) values for the rows whose checksum changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9   = "9d5f1879dcd93ff3b6990e3684a068b9"
    15  = "2d18f785f1b590a19e1788251be63ba2"
    17  = "e977eef88d24d0b433f0c0a866f0addb"
    89  = "3a425473b901d99eeb2f8f05d1a7a9da"
    99  = "7332e19db9d80de1248db805e60f9312"
    110 = "a0cab0e46f110ea81f706b2fc5953f20"
    121 = "20c9262e263ca6a7eba82d39434d7337"
    136 = "8ce95b65453f057e5cb46acb1c5b265e"
    159 = "6c56b0d081d323241d4926c5d3e90bda"
    169 = "d321d6ac78ad3d5397984693326c7178"
    175 = "6988a7d8079cd0690a573f7b9e2adedb"
    180 = "9ff250cc2296e8b04e2e9c55eb7b492a"
    183 = "06d85141b36781cfefde172ffae2ddd2"
    200 = "4e6c4a52cb99b938d83482275799e8ce"
    213 = "289d9c7f686850f0271f00b042591a5a"
    227 = "82760c335d1800fd1aeb50687d6f826e"
    232 = "3f0a589ba5292d038af5d7e15f995d2b"
    281 = "181895aa68478a8ce5e37e3a6123fdf6"
    342 = "100420a40e983f6c7d3d44f9d706e436"
    509 = "5fc2b40d76b385d2ea400730d0209360"
    511 = "b3c0471f6ab03fe79ed3515cd46b22cc"
    527 = "ef7389faaae9fb09633caaf053419f0b"
    559 = "a43aad2a42277be6fc85233bafe81f21"
    584 = "2dfdedb2c6659147cc3aefedac967c38"
    628 = "ed557febb7ee35ac12aced697d42ab84"
    639 = "0f0be0e84f78e6d2ed29c1de61ab4256"
    641 = "14ecbdeb7196c978887c9659fe12af43"
    667 = "e15e5c6bfedbf0f8a3460c19d2705282"
    677 = "16b63d480f3d50d78a869c19ab998727"
    692 = "4fc5fa4b3dd3ce2d2f863a4ac7f1255b"
    727 = "135f8f298b424ab3e3044b9d76c5be89"
    745 = "7efd05d52f32614d45b38870ab9df00a"
    780 = "7b32c2e2138ad20d6de90800ca768f42"
    831 = "fe026ef77817745f67b3c7b7b4ae0cd2"
    842 = "9497200fa07dee483d3cc5e671b64825"
    866 = "2cfba06a9d4752f41f6f3629d4306aaf"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
